$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.019565613109248
$ws.Range("D2").Value = 1.038363776616237
$ws.Range("E2").Value = 1.033524582211499
$ws.Range("F2").Value = 1.044693412768122
$ws.Range("I2").Value = 1.057164429023977
$ws.Range("J2").Value = 1.041176463476297
$ws.Range("K2").Value = 1.049313314935302
$ws.Range("L2").Value = 1.044535817776128
$ws.Range("M2").Value = 1.055563326922633
$ws.Range("N2").Value = 1.017546734616847
$ws.Range("P2").Value = 1.05254693794818
$ws.Range("C3").Value = 1.022896725595682
$ws.Range("D3").Value = 1.040729836616144
$ws.Range("E3").Value = 1.036083576252714
$ws.Range("F3").Value = 1.047168449425325
$ws.Range("I3").Value = 1.058098923253964
$ws.Range("J3").Value = 1.042782084564859
$ws.Range("K3").Value = 1.050866802538834
$ws.Range("L3").Value = 1.046274674888213
$ws.Range("M3").Value = 1.057231380148697
$ws.Range("N3").Value = 1.018085274038615
$ws.Range("P3").Value = 1.053867067492624
$ws.Range("C4").Value = 1.025017036395862
$ws.Range("D4").Value = 1.042238482024075
$ws.Range("E4").Value = 1.037717795090245
$ws.Range("F4").Value = 1.048749179832147
$ws.Range("I4").Value = 1.058684761312362
$ws.Range("J4").Value = 1.043801752675303
$ws.Range("K4").Value = 1.051852512658848
$ws.Range("L4").Value = 1.047381352715187
$ws.Range("M4").Value = 1.058292801170245
$ws.Range("N4").Value = 1.018427228921907
$ws.Range("P4").Value = 1.054707095392226
$ws.Range("C5").Value = 1.025900215219379
$ws.Range("D5").Value = 1.042867474758817
$ws.Range("E5").Value = 1.038399762519888
$ws.Range("F5").Value = 1.049408853207942
$ws.Range("I5").Value = 1.058926594808097
$ws.Range("J5").Value = 1.044225894207655
$ws.Range("K5").Value = 1.052262318628374
$ws.Range("L5").Value = 1.047842263686971
$ws.Range("M5").Value = 1.058734814382691
$ws.Range("N5").Value = 1.018569455093195
$ws.Range("P5").Value = 1.055056912412632
$ws.Range("C6").Value = 1.02604803004435
$ws.Range("D6").Value = 1.042972781146895
$ws.Range("E6").Value = 1.038513974599418
$ws.Range("F6").Value = 1.049519333076615
$ws.Range("I6").Value = 1.058966940245875
$ws.Range("J6").Value = 1.044296846469474
$ws.Range("K6").Value = 1.052330860220333
$ws.Range("L6").Value = 1.047919400867865
$ws.Range("M6").Value = 1.058808785911359
$ws.Range("N6").Value = 1.018593246479842
$ws.Range("P6").Value = 1.05511545476401
$ws.Range("C7").Value = 1.0250288694164
$ws.Range("D7").Value = 1.042246907108795
$ws.Range("E7").Value = 1.03772692729303
$ws.Range("F7").Value = 1.048758013394807
$ws.Range("I7").Value = 1.058688010123949
$ws.Range("J7").Value = 1.043807437748239
$ws.Range("K7").Value = 1.051858006421271
$ws.Range("L7").Value = 1.047387528354684
$ws.Range("M7").Value = 1.058298723800216
$ws.Range("N7").Value = 1.018429135336332
$ws.Range("P7").Value = 1.054711782666683
$ws.Range("C8").Value = 1.020698785307843
$ws.Range("D8").Value = 1.039168093390571
$ws.Range("E8").Value = 1.034393955414156
$ws.Range("F8").Value = 1.045534233808585
$ws.Range("I8").Value = 1.057484182802352
$ws.Range("J8").Value = 1.041723135769392
$ws.Range("K8").Value = 1.049842411222299
$ws.Range("L8").Value = 1.045127352562519
$ws.Range("M8").Value = 1.056130814343523
$ws.Range("N8").Value = 1.017730103128347
$ws.Range("P8").Value = 1.052996058714307
$ws.Range("C9").Value = 1.012789522716919
$ws.Range("D9").Value = 1.033566516806173
$ws.Range("E9").Value = 1.028349745568142
$ws.Range("F9").Value = 1.039689226024018
$ws.Range("I9").Value = 1.055215854562518
$ws.Range("J9").Value = 1.037898595184615
$ws.Range("K9").Value = 1.046137502843677
$ws.Range("L9").Value = 1.04099902335688
$ws.Range("M9").Value = 1.052169596812358
$ws.Range("N9").Value = 1.016447089158546
$ws.Range("P9").Value = 1.049861067226127
$ws.Range("C10").Value = 1.007360334551046
$ws.Range("D10").Value = 1.02975538096333
$ws.Range("E10").Value = 1.024252989676134
$ws.Range("F10").Value = 1.035759289862258
$ws.Range("I10").Value = 1.05363309498445
$ws.Range("J10").Value = 1.035284940215133
$ws.Range("K10").Value = 1.043607401522469
$ws.Range("L10").Value = 1.038197734006262
$ws.Range("M10").Value = 1.049511581434225
$ws.Range("N10").Value = 1.015574054659047
$ws.Range("P10").Value = 1.047808312724256
$ws.Range("C11").Value = 1.005325245380084
$ws.Range("D11").Value = 1.028466444149996
$ws.Range("E11").Value = 1.022892747719128
$ws.Range("F11").Value = 1.034700726176672
$ws.Range("I11").Value = 1.053180059532728
$ws.Range("J11").Value = 1.034478153785232
$ws.Range("K11").Value = 1.04287077672385
$ws.Range("L11").Value = 1.037395479153714
$ws.Range("M11").Value = 1.048996622308277
$ws.Range("N11").Value = 1.015338056203529
$ws.Range("P11").Value = 1.047832334509328
$ws.Range("C12").Value = 1.004698695008348
$ws.Range("D12").Value = 1.028128669544503
$ws.Range("E12").Value = 1.022547652185405
$ws.Range("F12").Value = 1.034551554798382
$ws.Range("I12").Value = 1.053103484526812
$ws.Range("J12").Value = 1.034305520312712
$ws.Range("K12").Value = 1.042736106847457
$ws.Range("L12").Value = 1.037255259093569
$ws.Range("M12").Value = 1.049045445587898
$ws.Range("N12").Value = 1.015307640235313
$ws.Range("P12").Value = 1.04819418539537
$ws.Range("C13").Value = 1.005116708645108
$ws.Range("D13").Value = 1.028503383308191
$ws.Range("E13").Value = 1.022963479534471
$ws.Range("F13").Value = 1.035102024140517
$ws.Range("I13").Value = 1.053318217628108
$ws.Range("J13").Value = 1.034614252616196
$ws.Range("K13").Value = 1.043061873913844
$ws.Range("L13").Value = 1.037620970926552
$ws.Range("M13").Value = 1.0495443846546
$ws.Range("N13").Value = 1.015435029304933
$ws.Range("P13").Value = 1.0488629867657
$ws.Range("C14").Value = 1.005882079575426
$ws.Range("D14").Value = 1.029087770776436
$ws.Range("E14").Value = 1.023598453855033
$ws.Range("F14").Value = 1.035803728788418
$ws.Range("I14").Value = 1.053599713091157
$ws.Range("J14").Value = 1.035047495579844
$ws.Range("K14").Value = 1.043497497926172
$ws.Range("L14").Value = 1.0381050538254
$ws.Range("M14").Value = 1.050096741930327
$ws.Range("N14").Value = 1.015595269328655
$ws.Range("P14").Value = 1.049471641031184
$ws.Range("C15").Value = 1.006298192807984
$ws.Range("D15").Value = 1.029389348252489
$ws.Range("E15").Value = 1.023923667862669
$ws.Range("F15").Value = 1.036136262090985
$ws.Range("I15").Value = 1.053734696050067
$ws.Range("J15").Value = 1.035261768363392
$ws.Range("K15").Value = 1.043708549033858
$ws.Range("L15").Value = 1.038338591481449
$ws.Range("M15").Value = 1.05033910069037
$ws.Range("N15").Value = 1.015670130638151
$ws.Range("P15").Value = 1.049700482234765
$ws.Range("C16").Value = 1.008514784710888
$ws.Range("D16").Value = 1.030933610932251
$ws.Range("E16").Value = 1.02557954062935
$ws.Range("F16").Value = 1.037715299437203
$ws.Range("I16").Value = 1.054380416918646
$ws.Range("J16").Value = 1.036318882205362
$ws.Range("K16").Value = 1.044730836263059
$ws.Range("L16").Value = 1.039466444788079
$ws.Range("M16").Value = 1.05140064150175
$ws.Range("N16").Value = 1.016018726031589
$ws.Range("P16").Value = 1.050501081512441
$ws.Range("C17").Value = 1.009815154429401
$ws.Range("D17").Value = 1.03181224191488
$ws.Range("E17").Value = 1.026517716253557
$ws.Range("F17").Value = 1.038555540093979
$ws.Range("I17").Value = 1.054725051489917
$ws.Range("J17").Value = 1.036901306737848
$ws.Range("K17").Value = 1.045284567708416
$ws.Range("L17").Value = 1.040076238072992
$ws.Range("M17").Value = 1.051919714049679
$ws.Range("N17").Value = 1.016201516600763
$ws.Range("P17").Value = 1.050783303420993
$ws.Range("C18").Value = 1.010443655950918
$ws.Range("D18").Value = 1.032187373117821
$ws.Range("E18").Value = 1.026910044966932
$ws.Range("F18").Value = 1.038808747914144
$ws.Range("I18").Value = 1.054835336148584
$ws.Range("J18").Value = 1.037117838813717
$ws.Range("K18").Value = 1.045472714563206
$ws.Range("L18").Value = 1.040279898450404
$ws.Range("M18").Value = 1.051989606959945
$ws.Range("N18").Value = 1.016253877182686
$ws.Range("P18").Value = 1.050602895944536
$ws.Range("C19").Value = 1.010458026810492
$ws.Range("D19").Value = 1.032100665927959
$ws.Range("E19").Value = 1.026801062092414
$ws.Range("F19").Value = 1.038523133050173
$ws.Range("I19").Value = 1.054731845924858
$ws.Range("J19").Value = 1.036999328885858
$ws.Range("K19").Value = 1.045325686945031
$ws.Range("L19").Value = 1.040110520652143
$ws.Range("M19").Value = 1.051647302646224
$ws.Range("N19").Value = 1.016188288908987
$ws.Range("P19").Value = 1.050009907029222
$ws.Range("C20").Value = 1.008766454141842
$ws.Range("D20").Value = 1.030741901463562
$ws.Range("E20").Value = 1.025312655087034
$ws.Range("F20").Value = 1.036777320467577
$ws.Range("I20").Value = 1.054047095253236
$ws.Range("J20").Value = 1.03596357670766
$ws.Range("K20").Value = 1.044264904413153
$ws.Range("L20").Value = 1.038924550568004
$ws.Range("M20").Value = 1.050202907409379
$ws.Range("N20").Value = 1.015800912139054
$ws.Range("P20").Value = 1.048344986493345
$ws.Range("C21").Value = 1.004629974151098
$ws.Range("D21").Value = 1.027818593418232
$ws.Range("E21").Value = 1.022170261268478
$ws.Range("F21").Value = 1.033715338991131
$ws.Range("I21").Value = 1.052797259221775
$ws.Range("J21").Value = 1.033935385434767
$ws.Range("K21").Value = 1.042291428549275
$ws.Range("L21").Value = 1.036743442133523
$ws.Range("M21").Value = 1.048084929649291
$ws.Range("N21").Value = 1.015117319598832
$ws.Range("P21").Value = 1.04662836127885
$ws.Range("C22").Value = 1.002005280048946
$ws.Range("D22").Value = 1.025976510704334
$ws.Range("E22").Value = 1.020194086123604
$ws.Range("F22").Value = 1.031805390618237
$ws.Range("I22").Value = 1.052004893242952
$ws.Range("J22").Value = 1.032657710833417
$ws.Range("K22").Value = 1.041050353959969
$ws.Range("L22").Value = 1.035375835536661
$ws.Range("M22").Value = 1.046772054131277
$ws.Range("N22").Value = 1.014688609085571
$ws.Range("P22").Value = 1.045589316335221
$ws.Range("C23").Value = 1.003401203619572
$ws.Range("D23").Value = 1.026955786901454
$ws.Range("E23").Value = 1.021244410021586
$ws.Range("F23").Value = 1.032820484763995
$ws.Range("I23").Value = 1.052427109293688
$ws.Range("J23").Value = 1.033337376838747
$ws.Range("K23").Value = 1.041710612459066
$ws.Range("L23").Value = 1.036103092959954
$ws.Range("M23").Value = 1.04747021622374
$ws.Range("N23").Value = 1.01491666389711
$ws.Range("P23").Value = 1.046141860747784
$ws.Range("C24").Value = 1.00879836394863
$ws.Range("D24").Value = 1.030750987464934
$ws.Range("E24").Value = 1.02532021053814
$ws.Range("F24").Value = 1.036760216990305
$ws.Range("I24").Value = 1.054041772517104
$ws.Range("J24").Value = 1.035961884626094
$ws.Range("K24").Value = 1.044258776426552
$ws.Range("L24").Value = 1.038916801522165
$ws.Range("M24").Value = 1.050171135555367
$ws.Range("N24").Value = 1.015797283185104
$ws.Range("P24").Value = 1.048279437649208
$ws.Range("C25").Value = 1.014868784302107
$ws.Range("D25").Value = 1.035036297538419
$ws.Range("E25").Value = 1.029933393424
$ws.Range("F25").Value = 1.041220515670911
$ws.Range("I25").Value = 1.055820093851709
$ws.Range("J25").Value = 1.038905910601357
$ws.Range("K25").Value = 1.047114020232099
$ws.Range("L25").Value = 1.042084146571857
$ws.Range("M25").Value = 1.053210944328785
$ws.Range("N25").Value = 1.01678504338056
$ws.Range("P25").Value = 1.050685212986644
